$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.099.07'
$ws.Range('E2').Value = '  +4.92%  '
$ws.Range('D3').Value = '2.617.58'
$ws.Range('E3').Value = '  +5.48%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '182.35'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.01%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.524'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.99%  '
$ws.Range('D9').Value = '2.617.42'
$ws.Range('E9').Value = '  +5.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.166'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +14.47%  '
$ws.Range('E12').Value = '  +4.55%  '
$ws.Range('E13').Value = '  +1.75%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.092.32'
$ws.Range('E14').Value = '  +5.36%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.84'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.15%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000184'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.93%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '70.915.44'
$ws.Range('E17').Value = '  +4.70%  '
$ws.Range('D18').Value = '2.626.50'
$ws.Range('E18').Value = '  +4.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '381.38'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.43%  '
$ws.Range('E20').Value = '  +6.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.20'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.23'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.43'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.73%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  +10.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.81'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +11.21%  '
$ws.Range('D28').Value = '2.746.77'
$ws.Range('E28').Value = '  +5.50%  '
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '0.0₃0953'
$ws.Range('E30').Value = '  +6.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '531.38'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.09'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.34'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.00%  '
$ws.Range('E34').Value = '  +4.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '164.04'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.23'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.94'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.64%  '
$ws.Range('E40').Value = '  +7.00%  '
$ws.Range('E41').Value = '  +6.09%  '
$ws.Range('E42').Value = '  +6.09%  '
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.59'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.332'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.03%  '
$ws.Range('E47').Value = '  +4.43%  '
$ws.Range('E48').Value = '  +4.13%  '
$ws.Range('D49').Value = '0.0₆0271'
$ws.Range('E49').Value = '  +6.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.535'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.72%  '
$ws.Range('E51').Value = '  +7.13%  '
